$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new journal entry on row 9
$ws.Range("B9").Value = 43902
$ws.Range("C9").Value = 0.6777777777777777
$ws.Range("F9").Value = "CLion"
$ws.Range("G9").Value = "Bataille Navale"
$ws.Range("H9").Value = "Programmation du jeu"
$ws.Range("I9").Value = "Création des différentes fonctions"

# Move the active selection to J2
$ws.Range("J2").Select()
